$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UQ")

# Update default glazing uncertainty: the Std Dev. column (G) for every
# StandardGlazing optical/thermal property row (rows 11-19) moves from
# 0.01 to 0.04.
$ws.Range("G11:G19").Value = 0.04

# Leave the sheet's selection where the author left it after making the edit.
$ws.Range("G20").Select()
